$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-28 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-29 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("993÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "175÷6=", 2) | Out-Null
$d.Content.Find.Execute("296÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "567÷9=", 2) | Out-Null
$d.Content.Find.Execute("310÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "588÷2=", 2) | Out-Null
$d.Content.Find.Execute("863÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "853÷2=", 2) | Out-Null
$d.Content.Find.Execute("390÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "349÷8=", 2) | Out-Null
$d.Content.Find.Execute("805÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "499÷7=", 2) | Out-Null
$d.Content.Find.Execute("295÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "283÷7=", 2) | Out-Null
$d.Content.Find.Execute("176÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "864÷5=", 2) | Out-Null
$d.Content.Find.Execute("672÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "387÷2=", 2) | Out-Null
$d.Content.Find.Execute("427÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "355÷7=", 2) | Out-Null
$d.Content.Find.Execute("725÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "599÷7=", 2) | Out-Null
$d.Content.Find.Execute("519÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "970÷2=", 2) | Out-Null
$d.Content.Find.Execute("789÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "385÷9=", 2) | Out-Null
$d.Content.Find.Execute("737÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "343÷7=", 2) | Out-Null
$d.Content.Find.Execute("550÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "496÷5=", 2) | Out-Null
$d.Content.Find.Execute("948÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "591÷5=", 2) | Out-Null
$d.Content.Find.Execute("920÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "289÷7=", 2) | Out-Null
$d.Content.Find.Execute("930÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "914÷8=", 2) | Out-Null
$d.Content.Find.Execute("576÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "509÷4=", 2) | Out-Null
$d.Content.Find.Execute("524÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "302÷3=", 2) | Out-Null
$d.Content.Find.Execute("684÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "192÷6=", 2) | Out-Null
$d.Content.Find.Execute("110÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "409÷7=", 2) | Out-Null
$d.Content.Find.Execute("854÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "644÷2=", 2) | Out-Null
$d.Content.Find.Execute("451÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "165÷3=", 2) | Out-Null
$d.Content.Find.Execute("158÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "143÷3=", 2) | Out-Null
